$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.201.98"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.926.66"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.51"
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.54"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000361"
$ws.Range("E11").Value = "  +6.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.90"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.69"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.552.75"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.93"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.920.00"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.16"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.341.10"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "448.31"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.81"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.40"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.68"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.55"
$ws.Range("E25").Value = "  +14.31%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.92"
$ws.Range("E26").Value = "  +13.85%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.62"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.16"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.52"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.132"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "689.63"
$ws.Range("E32").Value = "  -6.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.87"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("E34").Value = "  +20.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.91"
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.90"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.81"
$ws.Range("E37").Value = "  +7.77%  "
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.374"
$ws.Range("E41").Value = "  +10.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  +12.69%  "
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("E44").Value = "  +5.71%  "
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.17"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("E51").Value = "  -1.88%  "
